$d = $word.ActiveDocument

$d.Content.Find.Execute("12+16=28", $true, $false, $false, $false, $false, $true, 1, $false, "95-77=18", 2) | Out-Null
$d.Content.Find.Execute("20-17=3", $true, $false, $false, $false, $false, $true, 1, $false, "38-32=6", 2) | Out-Null
$d.Content.Find.Execute("72-48=24", $true, $false, $false, $false, $false, $true, 1, $false, "33+25=58", 2) | Out-Null
$d.Content.Find.Execute("53-30=23", $true, $false, $false, $false, $false, $true, 1, $false, "47-26=21", 2) | Out-Null
$d.Content.Find.Execute("11+62=73", $true, $false, $false, $false, $false, $true, 1, $false, "25+32=57", 2) | Out-Null
$d.Content.Find.Execute("35+32=67", $true, $false, $false, $false, $false, $true, 1, $false, "63-59=4", 2) | Out-Null
$d.Content.Find.Execute("40+50=90", $true, $false, $false, $false, $false, $true, 1, $false, "20+26=46", 2) | Out-Null
$d.Content.Find.Execute("23+61=84", $true, $false, $false, $false, $false, $true, 1, $false, "63-22=41", 2) | Out-Null
$d.Content.Find.Execute("55-25=30", $true, $false, $false, $false, $false, $true, 1, $false, "19+57=76", 2) | Out-Null
$d.Content.Find.Execute("91-57=34", $true, $false, $false, $false, $false, $true, 1, $false, "53+29=82", 2) | Out-Null
$d.Content.Find.Execute("69+26=95", $true, $false, $false, $false, $false, $true, 1, $false, "8+19=27", 2) | Out-Null
$d.Content.Find.Execute("67+30=97", $true, $false, $false, $false, $false, $true, 1, $false, "67+27=94", 2) | Out-Null
$d.Content.Find.Execute("19+49=68", $true, $false, $false, $false, $false, $true, 1, $false, "39-19=20", 2) | Out-Null
$d.Content.Find.Execute("95-94=1", $true, $false, $false, $false, $false, $true, 1, $false, "10+2=12", 2) | Out-Null
$d.Content.Find.Execute("67+31=98", $true, $false, $false, $false, $false, $true, 1, $false, "89-67=22", 2) | Out-Null
$d.Content.Find.Execute("19-0=19", $true, $false, $false, $false, $false, $true, 1, $false, "44-41=3", 2) | Out-Null
$d.Content.Find.Execute("82-32=50", $true, $false, $false, $false, $false, $true, 1, $false, "75-52=23", 2) | Out-Null
$d.Content.Find.Execute("29+39=68", $true, $false, $false, $false, $false, $true, 1, $false, "66-24=42", 2) | Out-Null
$d.Content.Find.Execute("95-36=59", $true, $false, $false, $false, $false, $true, 1, $false, "46-33=13", 2) | Out-Null
$d.Content.Find.Execute("76-36=40", $true, $false, $false, $false, $false, $true, 1, $false, "59-11=48", 2) | Out-Null
$d.Content.Find.Execute("36+28=64", $true, $false, $false, $false, $false, $true, 1, $false, "47-46=1", 2) | Out-Null
$d.Content.Find.Execute("92-71=21", $true, $false, $false, $false, $false, $true, 1, $false, "80-34=46", 2) | Out-Null
$d.Content.Find.Execute("92-10=82", $true, $false, $false, $false, $false, $true, 1, $false, "44+6=50", 2) | Out-Null
$d.Content.Find.Execute("25+56=81", $true, $false, $false, $false, $false, $true, 1, $false, "3+59=62", 2) | Out-Null
$d.Content.Find.Execute("33+14=47", $true, $false, $false, $false, $false, $true, 1, $false, "76-72=4", 2) | Out-Null
$d.Content.Find.Execute("46+48=94", $true, $false, $false, $false, $false, $true, 1, $false, "85-37=48", 2) | Out-Null
$d.Content.Find.Execute("29-15=14", $true, $false, $false, $false, $false, $true, 1, $false, "33-27=6", 2) | Out-Null
$d.Content.Find.Execute("43-19=24", $true, $false, $false, $false, $false, $true, 1, $false, "55+3=58", 2) | Out-Null
$d.Content.Find.Execute("53-53=0", $true, $false, $false, $false, $false, $true, 1, $false, "97-87=10", 2) | Out-Null
$d.Content.Find.Execute("7+33=40", $true, $false, $false, $false, $false, $true, 1, $false, "13+21=34", 2) | Out-Null
$d.Content.Find.Execute("13-12=1", $true, $false, $false, $false, $false, $true, 1, $false, "67-61=6", 2) | Out-Null
$d.Content.Find.Execute("11+17=28", $true, $false, $false, $false, $false, $true, 1, $false, "21-15=6", 2) | Out-Null
$d.Content.Find.Execute("21+60=81", $true, $false, $false, $false, $false, $true, 1, $false, "98-96=2", 2) | Out-Null
$d.Content.Find.Execute("67-59=8", $true, $false, $false, $false, $false, $true, 1, $false, "1+87=88", 2) | Out-Null
$d.Content.Find.Execute("19+51=70", $true, $false, $false, $false, $false, $true, 1, $false, "86-35=51", 2) | Out-Null
$d.Content.Find.Execute("80-78=2", $true, $false, $false, $false, $false, $true, 1, $false, "43-9=34", 2) | Out-Null
$d.Content.Find.Execute("78-7=71", $true, $false, $false, $false, $false, $true, 1, $false, "10+32=42", 2) | Out-Null
$d.Content.Find.Execute("84+9=93", $true, $false, $false, $false, $false, $true, 1, $false, "72-29=43", 2) | Out-Null
$d.Content.Find.Execute("62-52=10", $true, $false, $false, $false, $false, $true, 1, $false, "55-9=46", 2) | Out-Null
$d.Content.Find.Execute("73-44=29", $true, $false, $false, $false, $false, $true, 1, $false, "48+14=62", 2) | Out-Null
$d.Content.Find.Execute("79-31=48", $true, $false, $false, $false, $false, $true, 1, $false, "78-11=67", 2) | Out-Null
$d.Content.Find.Execute("34-9=25", $true, $false, $false, $false, $false, $true, 1, $false, "46+22=68", 2) | Out-Null
$d.Content.Find.Execute("44-22=22", $true, $false, $false, $false, $false, $true, 1, $false, "86-57=29", 2) | Out-Null
$d.Content.Find.Execute("38+48=86", $true, $false, $false, $false, $false, $true, 1, $false, "35-23=12", 2) | Out-Null
$d.Content.Find.Execute("72-35=37", $true, $false, $false, $false, $false, $true, 1, $false, "13+27=40", 2) | Out-Null
$d.Content.Find.Execute("53+37=90", $true, $false, $false, $false, $false, $true, 1, $false, "51-48=3", 2) | Out-Null
$d.Content.Find.Execute("83-15=68", $true, $false, $false, $false, $false, $true, 1, $false, "78+20=98", 2) | Out-Null
$d.Content.Find.Execute("96-63=33", $true, $false, $false, $false, $false, $true, 1, $false, "78-27=51", 2) | Out-Null
$d.Content.Find.Execute("92-51=41", $true, $false, $false, $false, $false, $true, 1, $false, "91-46=45", 2) | Out-Null
$d.Content.Find.Execute("37+0=37", $true, $false, $false, $false, $false, $true, 1, $false, "87-52=35", 2) | Out-Null
$d.Content.Find.Execute("44+37=81", $true, $false, $false, $false, $false, $true, 1, $false, "79+19=98", 2) | Out-Null
$d.Content.Find.Execute("43-30=13", $true, $false, $false, $false, $false, $true, 1, $false, "23-0=23", 2) | Out-Null
$d.Content.Find.Execute("67+4=71", $true, $false, $false, $false, $false, $true, 1, $false, "56-23=33", 2) | Out-Null
$d.Content.Find.Execute("90-77=13", $true, $false, $false, $false, $false, $true, 1, $false, "57+9=66", 2) | Out-Null
$d.Content.Find.Execute("32+52=84", $true, $false, $false, $false, $false, $true, 1, $false, "95-83=12", 2) | Out-Null
$d.Content.Find.Execute("57-37=20", $true, $false, $false, $false, $false, $true, 1, $false, "20-14=6", 2) | Out-Null
$d.Content.Find.Execute("43+12=55", $true, $false, $false, $false, $false, $true, 1, $false, "53-7=46", 2) | Out-Null
$d.Content.Find.Execute("50-17=33", $true, $false, $false, $false, $false, $true, 1, $false, "86-43=43", 2) | Out-Null
$d.Content.Find.Execute("90-69=21", $true, $false, $false, $false, $false, $true, 1, $false, "54-49=5", 2) | Out-Null
$d.Content.Find.Execute("55-49=6", $true, $false, $false, $false, $false, $true, 1, $false, "55-54=1", 2) | Out-Null
$d.Content.Find.Execute("63-46=17", $true, $false, $false, $false, $false, $true, 1, $false, "5+8=13", 2) | Out-Null
$d.Content.Find.Execute("24+58=82", $true, $false, $false, $false, $false, $true, 1, $false, "43-32=11", 2) | Out-Null
$d.Content.Find.Execute("41-10=31", $true, $false, $false, $false, $false, $true, 1, $false, "74+6=80", 2) | Out-Null
$d.Content.Find.Execute("55-42=13", $true, $false, $false, $false, $false, $true, 1, $false, "65-56=9", 2) | Out-Null
$d.Content.Find.Execute("26+41=67", $true, $false, $false, $false, $false, $true, 1, $false, "97-24=73", 2) | Out-Null
$d.Content.Find.Execute("6+58=64", $true, $false, $false, $false, $false, $true, 1, $false, "17+11=28", 2) | Out-Null
$d.Content.Find.Execute("16-10=6", $true, $false, $false, $false, $false, $true, 1, $false, "54-29=25", 2) | Out-Null
$d.Content.Find.Execute("44-33=11", $true, $false, $false, $false, $false, $true, 1, $false, "52+25=77", 2) | Out-Null
$d.Content.Find.Execute("80+0=80", $true, $false, $false, $false, $false, $true, 1, $false, "9+33=42", 2) | Out-Null
$d.Content.Find.Execute("24+59=83", $true, $false, $false, $false, $false, $true, 1, $false, "77-21=56", 2) | Out-Null
$d.Content.Find.Execute("12+43=55", $true, $false, $false, $false, $false, $true, 1, $false, "14+50=64", 2) | Out-Null
$d.Content.Find.Execute("12+66=78", $true, $false, $false, $false, $false, $true, 1, $false, "58-33=25", 2) | Out-Null
$d.Content.Find.Execute("2+4=6", $true, $false, $false, $false, $false, $true, 1, $false, "48-11=37", 2) | Out-Null
$d.Content.Find.Execute("36-31=5", $true, $false, $false, $false, $false, $true, 1, $false, "19+16=35", 2) | Out-Null
$d.Content.Find.Execute("89-3=86", $true, $false, $false, $false, $false, $true, 1, $false, "92-56=36", 2) | Out-Null
$d.Content.Find.Execute("4+82=86", $true, $false, $false, $false, $false, $true, 1, $false, "79-73=6", 2) | Out-Null
$d.Content.Find.Execute("96-38=58", $true, $false, $false, $false, $false, $true, 1, $false, "64-6=58", 2) | Out-Null
$d.Content.Find.Execute("91-87=4", $true, $false, $false, $false, $false, $true, 1, $false, "8+90=98", 2) | Out-Null
$d.Content.Find.Execute("75-39=36", $true, $false, $false, $false, $false, $true, 1, $false, "50+36=86", 2) | Out-Null
$d.Content.Find.Execute("78+5=83", $true, $false, $false, $false, $false, $true, 1, $false, "17+9=26", 2) | Out-Null
$d.Content.Find.Execute("8-7=1", $true, $false, $false, $false, $false, $true, 1, $false, "38+12=50", 2) | Out-Null
$d.Content.Find.Execute("74+14=88", $true, $false, $false, $false, $false, $true, 1, $false, "64-27=37", 2) | Out-Null
$d.Content.Find.Execute("15+42=57", $true, $false, $false, $false, $false, $true, 1, $false, "33+7=40", 2) | Out-Null
$d.Content.Find.Execute("31+35=66", $true, $false, $false, $false, $false, $true, 1, $false, "55-2=53", 2) | Out-Null
$d.Content.Find.Execute("90-41=49", $true, $false, $false, $false, $false, $true, 1, $false, "24+26=50", 2) | Out-Null
$d.Content.Find.Execute("71-65=6", $true, $false, $false, $false, $false, $true, 1, $false, "30+4=34", 2) | Out-Null
$d.Content.Find.Execute("33+29=62", $true, $false, $false, $false, $false, $true, 1, $false, "49-8=41", 2) | Out-Null
$d.Content.Find.Execute("21+42=63", $true, $false, $false, $false, $false, $true, 1, $false, "30-22=8", 2) | Out-Null
$d.Content.Find.Execute("45+54=99", $true, $false, $false, $false, $false, $true, 1, $false, "81+2=83", 2) | Out-Null
$d.Content.Find.Execute("36+24=60", $true, $false, $false, $false, $false, $true, 1, $false, "57-30=27", 2) | Out-Null
$d.Content.Find.Execute("22+65=87", $true, $false, $false, $false, $false, $true, 1, $false, "88-9=79", 2) | Out-Null
$d.Content.Find.Execute("92-62=30", $true, $false, $false, $false, $false, $true, 1, $false, "9+45=54", 2) | Out-Null
$d.Content.Find.Execute("54-52=2", $true, $false, $false, $false, $false, $true, 1, $false, "9+49=58", 2) | Out-Null
$d.Content.Find.Execute("79-17=62", $true, $false, $false, $false, $false, $true, 1, $false, "35+29=64", 2) | Out-Null
$d.Content.Find.Execute("32+20=52", $true, $false, $false, $false, $false, $true, 1, $false, "77-14=63", 2) | Out-Null
$d.Content.Find.Execute("97-61=36", $true, $false, $false, $false, $false, $true, 1, $false, "50+41=91", 2) | Out-Null
$d.Content.Find.Execute("47+33=80", $true, $false, $false, $false, $false, $true, 1, $false, "12+54=66", 2) | Out-Null
$d.Content.Find.Execute("7+92=99", $true, $false, $false, $false, $false, $true, 1, $false, "46+43=89", 2) | Out-Null
$d.Content.Find.Execute("34+55=89", $true, $false, $false, $false, $false, $true, 1, $false, "32-4=28", 2) | Out-Null
$d.Content.Find.Execute("12+64=76", $true, $false, $false, $false, $false, $true, 1, $false, "0+2=2", 2) | Out-Null
